{"js": "// Replace the date line and every addition/subtraction equation in the\n// practice-sheet table with the new values, matched by document order\n// (paragraph index) rather than by text search, because a few equations\n// (e.g. \"56-17=\") repeat with different intended replacements.\nconst OLD_VALUES = [\"2025-01-14 Tuesday\", \"19+19=\", \"62-53=\", \"27-17=\", \"53+11=\", \"18+0=\", \"47+16=\", \"39+60=\", \"58+28=\", \"9+72=\", \"18+79=\", \"36-35=\", \"56+18=\", \"93-78=\", \"33+0=\", \"29-16=\", \"39+19=\", \"21+68=\", \"71-26=\", \"85-10=\", \"8+67=\", \"94-52=\", \"6+75=\", \"31+19=\", \"56-17=\", \"57+22=\", \"76+7=\", \"99-54=\", \"13-2=\", \"56-17=\", \"91-0=\", \"59+11=\", \"28+35=\", \"64+14=\", \"12+61=\", \"38-36=\", \"85+13=\", \"11-10=\", \"30+60=\", \"99-69=\", \"4+25=\", \"4+28=\", \"97-37=\", \"51+34=\", \"94-76=\", \"82-79=\", \"30-10=\", \"57-27=\", \"28+25=\", \"94-5=\", \"73-55=\", \"10+2=\", \"19+80=\", \"1+72=\", \"64-32=\", \"41+23=\", \"58+40=\", \"31+35=\", \"0+20=\", \"54+37=\", \"62-25=\", \"19+49=\", \"79-16=\", \"49+16=\", \"83-46=\", \"68-61=\", \"94-83=\", \"89+7=\", \"88-64=\", \"39-11=\", \"95-73=\", \"88-72=\", \"43+32=\", \"61-37=\", \"36+16=\", \"34+20=\", \"40+51=\", \"42+30=\", \"16-9=\", \"4+29=\", \"25+32=\", \"81+9=\", \"87-70=\", \"53+4=\", \"16+82=\", \"38+15=\", \"56-49=\", \"66+27=\", \"24-23=\", \"45+37=\", \"91-32=\", \"87-8=\", \"96-37=\", \"48-34=\", \"93-18=\", \"92-21=\", \"12-6=\", \"66-11=\", \"68-24=\", \"11+27=\", \"90-45=\"];\nconst NEW_VALUES = [\"2025-01-15 Wednesday\", \"47+30=\", \"81-26=\", \"83-19=\", \"39+26=\", \"71-35=\", \"10+64=\", \"10+48=\", \"89-46=\", \"98-72=\", \"70+17=\", \"86+8=\", \"86-16=\", \"3+10=\", \"67+8=\", \"43+44=\", \"53+23=\", \"2+44=\", \"13+11=\", \"24+34=\", \"95-8=\", \"73-39=\", \"26+40=\", \"4+62=\", \"70-67=\", \"65+7=\", \"75-2=\", \"31+58=\", \"24+22=\", \"80-18=\", \"35+12=\", \"90-83=\", \"64-47=\", \"69-29=\", \"62-9=\", \"83+5=\", \"50-23=\", \"42+14=\", \"13+71=\", \"84-53=\", \"14+29=\", \"2+87=\", \"32+21=\", \"49+6=\", \"18+33=\", \"26-14=\", \"57+34=\", \"64-43=\", \"82-18=\", \"69-50=\", \"22+30=\", \"78+12=\", \"90-13=\", \"16+18=\", \"20+37=\", \"98-53=\", \"95-31=\", \"0+88=\", \"4+14=\", \"39-20=\", \"97-89=\", \"55-44=\", \"4+44=\", \"79+4=\", \"25+23=\", \"96-11=\", \"15+28=\", \"65+12=\", \"83-63=\", \"95-9=\", \"90-50=\", \"90-16=\", \"42-12=\", \"58-49=\", \"4-1=\", \"23-22=\", \"47-13=\", \"93-33=\", \"23-11=\", \"65+15=\", \"49+14=\", \"35-30=\", \"73-42=\", \"32-10=\", \"29-13=\", \"43-32=\", \"21+75=\", \"26+40=\", \"60-12=\", \"56+5=\", \"46+38=\", \"18+80=\", \"15+83=\", \"23+8=\", \"89-62=\", \"50-45=\", \"17+47=\", \"10+62=\", \"70-37=\", \"43-21=\", \"45-38=\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst n = Math.min(paragraphs.items.length, NEW_VALUES.length);\nfor (let i = 0; i < n; i++) {\n  const paragraph = paragraphs.items[i];\n  const currentText = paragraph.text;\n  const newText = NEW_VALUES[i];\n  if (currentText === newText) {\n    continue; // already correct, nothing to do\n  }\n  // currentText is expected to equal OLD_VALUES[i]; the replacement itself\n  // is always driven by position (not by matching OLD_VALUES) so duplicate\n  // old values (\"56-17=\" appears twice but maps to two different new\n  // values) are still handled correctly.\n  paragraph.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and every addition/subtraction equation in the\n# practice-sheet table with the new values, matched by document order\n# rather than by text search: a few equations (e.g. \"56-17=\") repeat\n# verbatim but map to two different replacements, so positional matching\n# is required for correctness.\n$NewValues = @(\n    \"2025-01-15 Wednesday\",\n    \"47+30=\",\n    \"81-26=\",\n    \"83-19=\",\n    \"39+26=\",\n    \"71-35=\",\n    \"10+64=\",\n    \"10+48=\",\n    \"89-46=\",\n    \"98-72=\",\n    \"70+17=\",\n    \"86+8=\",\n    \"86-16=\",\n    \"3+10=\",\n    \"67+8=\",\n    \"43+44=\",\n    \"53+23=\",\n    \"2+44=\",\n    \"13+11=\",\n    \"24+34=\",\n    \"95-8=\",\n    \"73-39=\",\n    \"26+40=\",\n    \"4+62=\",\n    \"70-67=\",\n    \"65+7=\",\n    \"75-2=\",\n    \"31+58=\",\n    \"24+22=\",\n    \"80-18=\",\n    \"35+12=\",\n    \"90-83=\",\n    \"64-47=\",\n    \"69-29=\",\n    \"62-9=\",\n    \"83+5=\",\n    \"50-23=\",\n    \"42+14=\",\n    \"13+71=\",\n    \"84-53=\",\n    \"14+29=\",\n    \"2+87=\",\n    \"32+21=\",\n    \"49+6=\",\n    \"18+33=\",\n    \"26-14=\",\n    \"57+34=\",\n    \"64-43=\",\n    \"82-18=\",\n    \"69-50=\",\n    \"22+30=\",\n    \"78+12=\",\n    \"90-13=\",\n    \"16+18=\",\n    \"20+37=\",\n    \"98-53=\",\n    \"95-31=\",\n    \"0+88=\",\n    \"4+14=\",\n    \"39-20=\",\n    \"97-89=\",\n    \"55-44=\",\n    \"4+44=\",\n    \"79+4=\",\n    \"25+23=\",\n    \"96-11=\",\n    \"15+28=\",\n    \"65+12=\",\n    \"83-63=\",\n    \"95-9=\",\n    \"90-50=\",\n    \"90-16=\",\n    \"42-12=\",\n    \"58-49=\",\n    \"4-1=\",\n    \"23-22=\",\n    \"47-13=\",\n    \"93-33=\",\n    \"23-11=\",\n    \"65+15=\",\n    \"49+14=\",\n    \"35-30=\",\n    \"73-42=\",\n    \"32-10=\",\n    \"29-13=\",\n    \"43-32=\",\n    \"21+75=\",\n    \"26+40=\",\n    \"60-12=\",\n    \"56+5=\",\n    \"46+38=\",\n    \"18+80=\",\n    \"15+83=\",\n    \"23+8=\",\n    \"89-62=\",\n    \"50-45=\",\n    \"17+47=\",\n    \"10+62=\",\n    \"70-37=\",\n    \"43-21=\",\n    \"45-38=\"\n)\n\n$d = $word.ActiveDocument\n$total = $d.Paragraphs.Count\n$idx = 0\n\nfor ($i = 1; $i -le $total; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    # Table rows carry an extra paragraph for the end-of-row mark\n    # (Range.Text == \"\\r\\a\"); skip those, they're not content cells.\n    $clean = $r.Text -replace \"[\\r\\x07]\", \"\"\n    if ($clean.Length -gt 0) {\n        if ($idx -lt $NewValues.Length) {\n            $newText = $NewValues[$idx]\n            if ($clean -ne $newText) {\n                $r.Text = $newText\n            }\n        }\n        $idx = $idx + 1\n    }\n}\n"}
